$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.694.86'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.948.47'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.01%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.12'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4840'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2939'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +1.43%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06813'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '112.42'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.42'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.35%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.949.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.02%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07650'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.511'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6907'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '296.37'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +7.79%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '30.718.94'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '13.32'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.692'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007701'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.202.09'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.570'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.26%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.800'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '168.41'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.37'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.94%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.181'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.82%  '
$ws.Range('E29').Value = '  +4.28%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.439'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.04%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.766'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +17.83%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.438'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.50%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05096'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7788'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.85%  '
$ws.Range('E35').Value = '  +2.72%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02076'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.732'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.702'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.45%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.046'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '111.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  +0.57%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8730'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.921'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '70.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.96%  '
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.373'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.469'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.96%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '48.64'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.05%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1253'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '35.63'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.88%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.2547'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.02%  '
